# Update "想去人数" (number of people interested) counts for a handful of
# conventions on both the "展览" sheet and the aggregated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions) sheet
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 5709
$wsExpo.Range("F5").Value = 959
$wsExpo.Range("F7").Value = 2594
$wsExpo.Range("F9").Value = 182
$wsExpo.Range("F11").Value = 96
$wsExpo.Range("F12").Value = 37
$wsExpo.Range("F13").Value = 2437
$wsExpo.Range("F14").Value = 475

# 全部类型 (All types) sheet - same events, shifted row numbers
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 5709
$wsAll.Range("F6").Value = 959
$wsAll.Range("F8").Value = 2594
$wsAll.Range("F10").Value = 182
$wsAll.Range("F13").Value = 96
$wsAll.Range("F14").Value = 37
$wsAll.Range("F15").Value = 2437
$wsAll.Range("F16").Value = 475
